{"js": "const replacements = [\n  [\"2024-09-26 Thursday\", \"2024-09-27 Friday\"],\n  [\"92-83=9\", \"60+12=72\"],\n  [\"34+48=82\", \"37+23=60\"],\n  [\"86-76=10\", \"68+2=70\"],\n  [\"69-12=57\", \"54+35=89\"],\n  [\"93-63=30\", \"62-51=11\"],\n  [\"76-39=37\", \"15+29=44\"],\n  [\"56-31=25\", \"85-19=66\"],\n  [\"52+17=69\", \"63-40=23\"],\n  [\"99-85=14\", \"48+15=63\"],\n  [\"30+2=32\", \"53+45=98\"],\n  [\"65-33=32\", \"58+27=85\"],\n  [\"34+60=94\", \"55+19=74\"],\n  [\"97-93=4\", \"7+55=62\"],\n  [\"26+45=71\", \"14+76=90\"],\n  [\"70+5=75\", \"65+16=81\"],\n  [\"93-16=77\", \"5+10=15\"],\n  [\"6+13=19\", \"24+59=83\"],\n  [\"58-57=1\", \"29-9=20\"],\n  [\"66-8=58\", \"91-86=5\"],\n  [\"53+2=55\", \"4+24=28\"],\n  [\"66-47=19\", \"99-2=97\"],\n  [\"25+49=74\", \"72-10=62\"],\n  [\"19+35=54\", \"25+37=62\"],\n  [\"58-2=56\", \"42-40=2\"],\n  [\"44-10=34\", \"80-79=1\"],\n  [\"27+16=43\", \"46+5=51\"],\n  [\"26-7=19\", \"83-16=67\"],\n  [\"57+22=79\", \"33+63=96\"],\n  [\"61+12=73\", \"3+63=66\"],\n  [\"25+55=80\", \"85-3=82\"],\n  [\"91-72=19\", \"65+11=76\"],\n  [\"69-56=13\", \"34-28=6\"],\n  [\"73+18=91\", \"7+63=70\"],\n  [\"77-23=54\", \"46-44=2\"],\n  [\"16+82=98\", \"1+75=76\"],\n  [\"55-28=27\", \"48-15=33\"],\n  [\"95-63=32\", \"96-49=47\"],\n  [\"11-3=8\", \"73+24=97\"],\n  [\"39+54=93\", \"61-53=8\"],\n  [\"22+11=33\", \"82-11=71\"],\n  [\"23+11=34\", \"36+17=53\"],\n  [\"75+16=91\", \"48-28=20\"],\n  [\"20+70=90\", \"54+18=72\"],\n  [\"37+33=70\", \"82-41=41\"],\n  [\"74+9=83\", \"0+46=46\"],\n  [\"1+88=89\", \"97-61=36\"],\n  [\"59+23=82\", \"79+6=85\"],\n  [\"4+56=60\", \"52-47=5\"],\n  [\"32-12=20\", \"45+32=77\"],\n  [\"75-71=4\", \"65+28=93\"],\n  [\"83-76=7\", \"68-11=57\"],\n  [\"40+42=82\", \"56-12=44\"],\n  [\"96-4=92\", \"2+77=79\"],\n  [\"90-71=19\", \"89-34=55\"],\n  [\"94+2=96\", \"46+0=46\"],\n  [\"1+17=18\", \"12+30=42\"],\n  [\"24+23=47\", \"51-40=11\"],\n  [\"5+7=12\", \"95-48=47\"],\n  [\"82-20=62\", \"25+29=54\"],\n  [\"60-54=6\", \"42-41=1\"],\n  [\"23-16=7\", \"37+19=56\"],\n  [\"64-38=26\", \"24+24=48\"],\n  [\"9+15=24\", \"87-20=67\"],\n  [\"92-75=17\", \"83-53=30\"],\n  [\"79+19=98\", \"91-54=37\"],\n  [\"67-56=11\", \"85-59=26\"],\n  [\"2+66=68\", \"36+52=88\"],\n  [\"62-3=59\", \"88-65=23\"],\n  [\"30+49=79\", \"78-42=36\"],\n  [\"75-61=14\", \"59-31=28\"],\n  [\"27+31=58\", \"12-11=1\"],\n  [\"56-49=7\", \"89-11=78\"],\n  [\"48+6=54\", \"74-8=66\"],\n  [\"46-7=39\", \"14+57=71\"],\n  [\"52+14=66\", \"61-18=43\"],\n  [\"64-57=7\", \"63+3=66\"],\n  [\"28+32=60\", \"53+18=71\"],\n  [\"70+19=89\", \"4+39=43\"],\n  [\"4+93=97\", \"68-42=26\"],\n  [\"79-2=77\", \"36+5=41\"],\n  [\"32+67=99\", \"49+9=58\"],\n  [\"8+58=66\", \"37+22=59\"],\n  [\"1+87=88\", \"49-48=1\"],\n  [\"58-30=28\", \"43+35=78\"],\n  [\"76-5=71\", \"39+55=94\"],\n  [\"11+82=93\", \"5+51=56\"],\n  [\"78-52=26\", \"49-33=16\"],\n  [\"80-73=7\", \"10+59=69\"],\n  [\"46+49=95\", \"16+45=61\"],\n  [\"52+13=65\", \"11+80=91\"],\n  [\"49-44=5\", \"79-64=15\"],\n  [\"33-6=27\", \"73+12=85\"],\n  [\"20+15=35\", \"24-17=7\"],\n  [\"49+38=87\", \"39-32=7\"],\n  [\"58-29=29\", \"70-13=57\"],\n  [\"19-13=6\", \"73-69=4\"],\n  [\"66+3=69\", \"0+54=54\"],\n  [\"35+14=49\", \"46+44=90\"],\n  [\"41+40=81\", \"45-16=29\"],\n  [\"97-41=56\", \"81-80=1\"],\n];\n\nconst body = context.document.body;\nfor (const [before, after] of replacements) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: true });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + before);\n  }\n  for (const item of results.items) {\n    item.insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Replace-Text($before, $after) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $before\n    $find.Replacement.Text = $after\n    # wdFindContinue=1 wrap, wdReplaceAll=2 -> replaces every (exact, case-sensitive) match of $before\n    $found = $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $after, 2)\n    if (-not $found) {\n        throw \"Replace-Text: no match found for '$before'\"\n    }\n}\n\nReplace-Text '2024-09-26 Thursday' '2024-09-27 Friday'\nReplace-Text '92-83=9' '60+12=72'\nReplace-Text '34+48=82' '37+23=60'\nReplace-Text '86-76=10' '68+2=70'\nReplace-Text '69-12=57' '54+35=89'\nReplace-Text '93-63=30' '62-51=11'\nReplace-Text '76-39=37' '15+29=44'\nReplace-Text '56-31=25' '85-19=66'\nReplace-Text '52+17=69' '63-40=23'\nReplace-Text '99-85=14' '48+15=63'\nReplace-Text '30+2=32' '53+45=98'\nReplace-Text '65-33=32' '58+27=85'\nReplace-Text '34+60=94' '55+19=74'\nReplace-Text '97-93=4' '7+55=62'\nReplace-Text '26+45=71' '14+76=90'\nReplace-Text '70+5=75' '65+16=81'\nReplace-Text '93-16=77' '5+10=15'\nReplace-Text '6+13=19' '24+59=83'\nReplace-Text '58-57=1' '29-9=20'\nReplace-Text '66-8=58' '91-86=5'\nReplace-Text '53+2=55' '4+24=28'\nReplace-Text '66-47=19' '99-2=97'\nReplace-Text '25+49=74' '72-10=62'\nReplace-Text '19+35=54' '25+37=62'\nReplace-Text '58-2=56' '42-40=2'\nReplace-Text '44-10=34' '80-79=1'\nReplace-Text '27+16=43' '46+5=51'\nReplace-Text '26-7=19' '83-16=67'\nReplace-Text '57+22=79' '33+63=96'\nReplace-Text '61+12=73' '3+63=66'\nReplace-Text '25+55=80' '85-3=82'\nReplace-Text '91-72=19' '65+11=76'\nReplace-Text '69-56=13' '34-28=6'\nReplace-Text '73+18=91' '7+63=70'\nReplace-Text '77-23=54' '46-44=2'\nReplace-Text '16+82=98' '1+75=76'\nReplace-Text '55-28=27' '48-15=33'\nReplace-Text '95-63=32' '96-49=47'\nReplace-Text '11-3=8' '73+24=97'\nReplace-Text '39+54=93' '61-53=8'\nReplace-Text '22+11=33' '82-11=71'\nReplace-Text '23+11=34' '36+17=53'\nReplace-Text '75+16=91' '48-28=20'\nReplace-Text '20+70=90' '54+18=72'\nReplace-Text '37+33=70' '82-41=41'\nReplace-Text '74+9=83' '0+46=46'\nReplace-Text '1+88=89' '97-61=36'\nReplace-Text '59+23=82' '79+6=85'\nReplace-Text '4+56=60' '52-47=5'\nReplace-Text '32-12=20' '45+32=77'\nReplace-Text '75-71=4' '65+28=93'\nReplace-Text '83-76=7' '68-11=57'\nReplace-Text '40+42=82' '56-12=44'\nReplace-Text '96-4=92' '2+77=79'\nReplace-Text '90-71=19' '89-34=55'\nReplace-Text '94+2=96' '46+0=46'\nReplace-Text '1+17=18' '12+30=42'\nReplace-Text '24+23=47' '51-40=11'\nReplace-Text '5+7=12' '95-48=47'\nReplace-Text '82-20=62' '25+29=54'\nReplace-Text '60-54=6' '42-41=1'\nReplace-Text '23-16=7' '37+19=56'\nReplace-Text '64-38=26' '24+24=48'\nReplace-Text '9+15=24' '87-20=67'\nReplace-Text '92-75=17' '83-53=30'\nReplace-Text '79+19=98' '91-54=37'\nReplace-Text '67-56=11' '85-59=26'\nReplace-Text '2+66=68' '36+52=88'\nReplace-Text '62-3=59' '88-65=23'\nReplace-Text '30+49=79' '78-42=36'\nReplace-Text '75-61=14' '59-31=28'\nReplace-Text '27+31=58' '12-11=1'\nReplace-Text '56-49=7' '89-11=78'\nReplace-Text '48+6=54' '74-8=66'\nReplace-Text '46-7=39' '14+57=71'\nReplace-Text '52+14=66' '61-18=43'\nReplace-Text '64-57=7' '63+3=66'\nReplace-Text '28+32=60' '53+18=71'\nReplace-Text '70+19=89' '4+39=43'\nReplace-Text '4+93=97' '68-42=26'\nReplace-Text '79-2=77' '36+5=41'\nReplace-Text '32+67=99' '49+9=58'\nReplace-Text '8+58=66' '37+22=59'\nReplace-Text '1+87=88' '49-48=1'\nReplace-Text '58-30=28' '43+35=78'\nReplace-Text '76-5=71' '39+55=94'\nReplace-Text '11+82=93' '5+51=56'\nReplace-Text '78-52=26' '49-33=16'\nReplace-Text '80-73=7' '10+59=69'\nReplace-Text '46+49=95' '16+45=61'\nReplace-Text '52+13=65' '11+80=91'\nReplace-Text '49-44=5' '79-64=15'\nReplace-Text '33-6=27' '73+12=85'\nReplace-Text '20+15=35' '24-17=7'\nReplace-Text '49+38=87' '39-32=7'\nReplace-Text '58-29=29' '70-13=57'\nReplace-Text '19-13=6' '73-69=4'\nReplace-Text '66+3=69' '0+54=54'\nReplace-Text '35+14=49' '46+44=90'\nReplace-Text '41+40=81' '45-16=29'\nReplace-Text '97-41=56' '81-80=1'\n"}
